# Updates the cryptos worksheet cell values to match the latest scrape.
# Values that look numeric (e.g. price column) must be forced to text so
# that Excel does not silently convert them to numbers (losing formatting
# such as trailing zeros or thousand-group separators used as literal dots).

function Set-TextValue {
    param($Worksheet, $Address, $Text)
    $Worksheet.Range($Address).NumberFormat = "@"
    $Worksheet.Range($Address).Value = $Text
    $Worksheet.Range($Address).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.420.04"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  -0.64%  "

Set-TextValue $ws "D5" "1.003"
$ws.Range("E5").Value = "  -0.53%  "

Set-TextValue $ws "D6" "305.43"
$ws.Range("E6").Value = "  -1.10%  "

Set-TextValue $ws "D7" "0.4501"
$ws.Range("E7").Value = "  -0.85%  "

Set-TextValue $ws "D8" "0.3583"
$ws.Range("E8").Value = "  -2.13%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D9" "46.28"
$ws.Range("E9").Value = "  +2.64%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws "D10" "0.07052"
$ws.Range("E10").Value = "  -1.00%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws "D11" "0.8887"
$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D12" "0.07773"
$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws "D13" "19.30"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.834.15"
$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D15" "5.268"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D16" "6.301"
$ws.Range("E16").Value = "  -0.74%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws "D17" "84.87"
$ws.Range("E17").Value = "  -1.47%  "

$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws "D18" "1.006"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D19" "0.000008513"
$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D20" "1.003"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.463.52"
$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws "D22" "14.16"
$ws.Range("E22").Value = "  -0.54%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws "D23" "4.955"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.043.67"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D25" "10.50"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D26" "1.953"
$ws.Range("E26").Value = "  -1.28%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D27" "150.97"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D28" "17.78"
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D29" "2.049"
$ws.Range("E29").Value = "  +2.79%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D30" "112.08"
$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D31" "4.832"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D32" "0.08678"
$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws "D33" "3.137"
$ws.Range("E33").Value = "  +3.45%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D34" "0.7407"
$ws.Range("E34").Value = "  +1.99%  "

Set-TextValue $ws "D35" "2.746"
$ws.Range("E35").Value = "  +8.17%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D36" "4.430"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D37" "1.106"
$ws.Range("E37").Value = "  -0.56%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D38" "1.066"
$ws.Range("E38").Value = "  -1.29%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D39" "0.01925"
$ws.Range("E39").Value = "  -0.14%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D40" "0.05112"
$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D41" "2.892"
$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws "D42" "0.5061"
$ws.Range("E42").Value = "  +1.34%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D43" "6.754"
$ws.Range("E43").Value = "  -2.74%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D44" "0.1504"
$ws.Range("E44").Value = "  -3.92%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D45" "8.041"
$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws "D46" "0.4687"
$ws.Range("E46").Value = "  +1.61%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws "D47" "1.003"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D48" "9.990"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws "D49" "99.91"
$ws.Range("E49").Value = "  -2.06%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D50" "1.575"
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D51" "0.05984"
$ws.Range("E51").Value = "  -0.23%  "
